$d = $word.ActiveDocument

# 1) "Juan Martin Santos Ayala - 202013610" -> replace hyphen separator with en dash
$r1 = $d.Paragraphs.Item(1).Range
$r1.Find.Execute(" - ", $false, $false, $false, $false, $false, $true, 1, $false, " – ", 2)

# 2) "Daniel Esteban Aguilera Figueroa – " -> append the second student id after the dash
$r2 = $d.Paragraphs.Item(2).Range
$r2.Find.Execute("– ", $false, $false, $false, $false, $false, $true, 1, $false, "– 202010592", 2)

# 3) Mark the following empty paragraph's mark as English (US)
$p3 = $d.Paragraphs.Item(3)
$p3.Range.LanguageID = "en-US"

# 4) Fix typo "Identificque" -> "Identifique"
$d.Content.Find.Execute("Identificque", $false, $false, $false, $false, $false, $true, 1, $false, "Identifique", 2)

# 5) Insert " in incluir" into the PDF conclusions sentence
$d.Content.Find.Execute("visualizaciones e imágenes.", $false, $false, $false, $false, $false, $true, 1, $false, "visualizaciones e in incluir imágenes.", 2)

# 6) Insert three new paragraphs (blank, new sentence, blank) right before the final blank paragraph
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$rLast = $pLast.Range
$rLast.InsertParagraphBefore()
$rLast.InsertParagraphBefore()
$rLast.InsertParagraphBefore()

$pNewText = $d.Paragraphs.Item($d.Paragraphs.Count - 2)
$pNewText.Range.InsertBefore("Basándonos en el framwork de visulización de Tamara, ")
